# Update the date heading, then every "AxB=" practice problem in the
# multiplication table to the new values from the target revision.
#
# Note on ordering: one new value ("344×5=", produced by turning
# "789×6=" into "344×5=") happens to equal another cell's *original*
# text (which itself must become "964×7="). To avoid the freshly
# written "344×5=" being caught by that later search-and-replace, the
# "344×5=" -> "964×7=" replacement is performed first, and the
# "789×6=" -> "344×5=" replacement is performed last.

$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-04-30 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-01 Thursday", 2) | Out-Null
$d.Content.Find.Execute("165×7=", $true, $false, $false, $false, $false, $true, 1, $false, "628×6=", 2) | Out-Null
$d.Content.Find.Execute("891×6=", $true, $false, $false, $false, $false, $true, 1, $false, "637×2=", 2) | Out-Null
$d.Content.Find.Execute("493×6=", $true, $false, $false, $false, $false, $true, 1, $false, "431×3=", 2) | Out-Null
$d.Content.Find.Execute("566×9=", $true, $false, $false, $false, $false, $true, 1, $false, "235×2=", 2) | Out-Null
$d.Content.Find.Execute("424×6=", $true, $false, $false, $false, $false, $true, 1, $false, "979×3=", 2) | Out-Null
$d.Content.Find.Execute("341×4=", $true, $false, $false, $false, $false, $true, 1, $false, "778×3=", 2) | Out-Null
$d.Content.Find.Execute("923×7=", $true, $false, $false, $false, $false, $true, 1, $false, "781×7=", 2) | Out-Null
$d.Content.Find.Execute("997×4=", $true, $false, $false, $false, $false, $true, 1, $false, "653×8=", 2) | Out-Null
$d.Content.Find.Execute("101×9=", $true, $false, $false, $false, $false, $true, 1, $false, "595×8=", 2) | Out-Null
$d.Content.Find.Execute("226×4=", $true, $false, $false, $false, $false, $true, 1, $false, "950×2=", 2) | Out-Null
$d.Content.Find.Execute("446×6=", $true, $false, $false, $false, $false, $true, 1, $false, "432×8=", 2) | Out-Null
$d.Content.Find.Execute("790×4=", $true, $false, $false, $false, $false, $true, 1, $false, "186×8=", 2) | Out-Null
$d.Content.Find.Execute("593×7=", $true, $false, $false, $false, $false, $true, 1, $false, "855×7=", 2) | Out-Null
$d.Content.Find.Execute("224×2=", $true, $false, $false, $false, $false, $true, 1, $false, "399×8=", 2) | Out-Null
$d.Content.Find.Execute("633×9=", $true, $false, $false, $false, $false, $true, 1, $false, "342×5=", 2) | Out-Null
$d.Content.Find.Execute("875×5=", $true, $false, $false, $false, $false, $true, 1, $false, "806×7=", 2) | Out-Null
$d.Content.Find.Execute("435×9=", $true, $false, $false, $false, $false, $true, 1, $false, "591×6=", 2) | Out-Null
$d.Content.Find.Execute("740×9=", $true, $false, $false, $false, $false, $true, 1, $false, "183×3=", 2) | Out-Null
$d.Content.Find.Execute("254×7=", $true, $false, $false, $false, $false, $true, 1, $false, "668×4=", 2) | Out-Null
$d.Content.Find.Execute("849×7=", $true, $false, $false, $false, $false, $true, 1, $false, "173×6=", 2) | Out-Null
$d.Content.Find.Execute("170×3=", $true, $false, $false, $false, $false, $true, 1, $false, "377×4=", 2) | Out-Null
$d.Content.Find.Execute("162×3=", $true, $false, $false, $false, $false, $true, 1, $false, "859×9=", 2) | Out-Null
$d.Content.Find.Execute("448×9=", $true, $false, $false, $false, $false, $true, 1, $false, "925×3=", 2) | Out-Null
$d.Content.Find.Execute("344×5=", $true, $false, $false, $false, $false, $true, 1, $false, "964×7=", 2) | Out-Null
$d.Content.Find.Execute("789×6=", $true, $false, $false, $false, $false, $true, 1, $false, "344×5=", 2) | Out-Null
